$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function SetDE($row, $dval, $eval) {
    $arr = New-Object 'object[,]' 1,2
    $arr[0,0] = $dval
    $arr[0,1] = $eval
    $ws.Range("D${row}:E${row}").Value2 = $arr
}

# Insert two new columns before column D (new quarter columns),
# shifting existing D:K data to F:M.
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formatting from column F into the newly inserted D:E
# columns for each of the three data blocks (Income Statement,
# Balance Sheet, Cash Flow Statement), so the new columns pick up the
# same date/number formats as the rest of the table.
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns (D, E) with the latest period
# data, and correct a handful of restated historical figures that
# moved along with the column shift.
SetDE 7 43496 43404
SetDE 8 232300 224700
$ws.Range("H8").Value2 = 186000
$ws.Range("I8").Value2 = 177000
SetDE 9 64500 61400
SetDE 10 167800 163300
$ws.Range("H10").Value2 = 127100
$ws.Range("I10").Value2 = 123800
SetDE 12 42800 40000
$ws.Range("H12").Value2 = 69100
$ws.Range("I12").Value2 = 63300
SetDE 13 0 0
SetDE 14 0 0
SetDE 15 0 0
SetDE 17 169300 161600
$ws.Range("H17").Value2 = 147500
SetDE 18 63000 63100
$ws.Range("H18").Value2 = 38500
$ws.Range("I18").Value2 = 42500
SetDE 20 5700 4600
$ws.Range("H20").Value2 = 3000
SetDE 21 72300 71100
$ws.Range("H21").Value2 = 45200
$ws.Range("I21").Value2 = 47500
SetDE 22 0 0
SetDE 23 68700 67700
$ws.Range("H23").Value2 = 41500
$ws.Range("I23").Value2 = 43900
SetDE 24 -2500 3600
$ws.Range("H24").Value2 = 900
$ws.Range("I24").Value2 = 8900
SetDE 25 0 0
SetDE 26 71200 64100
$ws.Range("H26").Value2 = 40700
$ws.Range("I26").Value2 = 34900
SetDE 27 71200 64100
$ws.Range("H27").Value2 = 40700
$ws.Range("I27").Value2 = 34900
SetDE 28 0 0
SetDE 29 0 0
SetDE 30 0 0
SetDE 31 0 0
SetDE 32 -5700 -4600
$ws.Range("H32").Value2 = -3000
SetDE 33 71200 64100
$ws.Range("H33").Value2 = 40700
$ws.Range("I33").Value2 = 34900
SetDE 34 0 0
SetDE 35 71200 64100
$ws.Range("H35").Value2 = 40700
$ws.Range("I35").Value2 = 34900
SetDE 38 43496 43404
SetDE 41 551000 467600
SetDE 42 539200 584300
SetDE 43 321600 110000
SetDE 44 0 0
SetDE 45 21700 21200
SetDE 46 1433400 1183200
SetDE 47 0 0
SetDE 48 55000 53600
SetDE 49 120300 122000
SetDE 50 0 0
SetDE 51 0 0
SetDE 52 45100 38500
SetDE 53 0 0
SetDE 54 1653800 1397300
SetDE 57 9100 10900
SetDE 58 0 0
SetDE 59 391900 225900
SetDE 60 401000 236700
SetDE 61 0 0
SetDE 62 15000 22300
SetDE 63 0 0
SetDE 64 0 0
SetDE 65 0 0
SetDE 66 416000 259000
SetDE 68 0 0
SetDE 69 0 0
SetDE 70 0 0
SetDE 71 0 0
SetDE 72 619200 548000
SetDE 73 0 0
SetDE 74 0 0
SetDE 75 0 0
SetDE 76 1237700 1138300
SetDE 77 0 0
SetDE 80 43496 43404
SetDE 81 71200 64100
$ws.Range("H81").Value2 = 40700
$ws.Range("I81").Value2 = 34900
SetDE 83 3600 3400
SetDE 84 0 0
SetDE 85 0 0
SetDE 86 0 0
SetDE 87 0 0
SetDE 88 0 0
SetDE 89 31900 41600
SetDE 91 -2900 -4200
$ws.Range("F91").Value2 = -700
$ws.Range("G91").Value2 = -700
$ws.Range("H91").Value2 = -1500
$ws.Range("I91").Value2 = -1600
$ws.Range("J91").Value2 = -2500
SetDE 92 0 0
SetDE 93 0 0
SetDE 94 43800 -89400
$ws.Range("I94").Value2 = -134400
SetDE 96 0 0
SetDE 97 0 0
SetDE 98 0 0
SetDE 99 0 0
SetDE 100 6200 4900
SetDE 101 1500 -1200
SetDE 102 83300 -44100
$ws.Range("I102").Value2 = -98400

Write-Output "Quarterly financials update applied."
